$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.1506834993901052
$arr[1,0] = 0.1506001191993676
$arr[2,0] = 0.1506033809586
$arr[3,0] = 0.1506184217764712
$arr[4,0] = 0.1506217479726928
$arr[5,0] = 0.1506035282615699
$arr[6,0] = 0.1506434575360345
$arr[7,0] = 0.1511532265322515
$arr[8,0] = 0.1517901743450381
$arr[9,0] = 0.1521368116117472
$arr[10,0] = 0.1522762418229675
$arr[11,0] = 0.1522458500570281
$arr[12,0] = 0.1521481190352105
$arr[13,0] = 0.1520893190332231
$arr[14,0] = 0.1517686642803397
$arr[15,0] = 0.1515865128976159
$arr[16,0] = 0.1514870994838873
$arr[17,0] = 0.1514543600112503
$arr[18,0] = 0.1516053490774141
$arr[19,0] = 0.1521766034877317
$arr[20,0] = 0.1525975438852925
$arr[21,0] = 0.1523685298813149
$arr[22,0] = 0.1515968167070838
$arr[23,0] = 0.1509691936213926
$ws.Range("C2:C25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.1005736726637849
$arr[1,0] = 0.1007978339275439
$arr[2,0] = 0.100960746591749
$arr[3,0] = 0.1010335045919994
$arr[4,0] = 0.1010459711264282
$arr[5,0] = 0.1009617020222358
$arr[6,0] = 0.1006457254056166
$arr[7,0] = 0.1002260900633232
$arr[8,0] = 0.100039000092309
$arr[9,0] = 0.09998006562749495
$arr[10,0] = 0.09996150009044413
$arr[11,0] = 0.09996533184503775
$arr[12,0] = 0.09997846310781
$arr[13,0] = 0.09998699460826543
$arr[14,0] = 0.100043377243221
$arr[15,0] = 0.1000846612311186
$arr[16,0] = 0.1001108709945164
$arr[17,0] = 0.1001201687867272
$arr[18,0] = 0.100080011521662
$arr[19,0] = 0.09997450441448663
$arr[20,0] = 0.09992741150920281
$arr[21,0] = 0.09995054944541693
$arr[22,0] = 0.100082105944491
$arr[23,0] = 0.1003182753274245
$ws.Range("D2:D25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.1396915254104059
$arr[1,0] = 0.1396102605363012
$arr[2,0] = 0.1396130617121685
$arr[3,0] = 0.1396274676000235
$arr[4,0] = 0.1396306612126121
$arr[5,0] = 0.1396132022703185
$arr[6,0] = 0.1396525751937041
$arr[7,0] = 0.1401475027033179
$arr[8,0] = 0.1407654421211291
$arr[9,0] = 0.1411017313729168
$arr[10,0] = 0.1412370034830523
$arr[11,0] = 0.1412075178115657
$arr[12,0] = 0.1411127014784341
$arr[13,0] = 0.1410556557639175
$arr[14,0] = 0.1407445745287284
$arr[15,0] = 0.1405678641227261
$arr[16,0] = 0.1404714198056887
$arr[17,0] = 0.1404396579129354
$arr[18,0] = 0.1405861376677429
$arr[19,0] = 0.1411403362758747
$arr[20,0] = 0.1415487359093213
$arr[21,0] = 0.1413265410027016
$arr[22,0] = 0.1405778601597945
$arr[23,0] = 0.1399689281676899
$ws.Range("E2:E25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 2.447612376704498
$arr[1,0] = 2.411058583748812
$arr[2,0] = 2.389927942158664
$arr[3,0] = 2.381646935565811
$arr[4,0] = 2.380291794617818
$arr[5,0] = 2.389814926386919
$arr[6,0] = 2.434735744743548
$arr[7,0] = 2.533273707637761
$arr[8,0] = 2.612086996391866
$arr[9,0] = 2.649345923571275
$arr[10,0] = 2.663657843630489
$arr[11,0] = 2.660566488450826
$arr[12,0] = 2.65051930683893
$arr[13,0] = 2.644391541257505
$arr[14,0] = 2.609680383338599
$arr[15,0] = 2.588746795704679
$arr[16,0] = 2.57683869061799
$arr[17,0] = 2.572829524297447
$arr[18,0] = 2.59096151124956
$arr[19,0] = 2.653464900058424
$arr[20,0] = 2.695496654110883
$arr[21,0] = 2.672955164908416
$arr[22,0] = 2.589959843127247
$arr[23,0] = 2.505493110964139
$ws.Range("F2:F25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.002499528841394606
$arr[1,0] = 0.002504394110527009
$arr[2,0] = 0.002507542966394094
$arr[3,0] = 0.002508866912270074
$arr[4,0] = 0.002509089218164611
$arr[5,0] = 0.002507560656470031
$arr[6,0] = 0.002501172935996968
$arr[7,0] = 0.002489922381588402
$arr[8,0] = 0.002482425708846821
$arr[9,0] = 0.002479180446426241
$arr[10,0] = 0.002477975138054288
$arr[11,0] = 0.002478233675257502
$arr[12,0] = 0.002479080812708475
$arr[13,0] = 0.002479602778810161
$arr[14,0] = 0.002482641103866756
$arr[15,0] = 0.002484547191223119
$arr[16,0] = 0.002485659061408373
$arr[17,0] = 0.002486038194566389
$arr[18,0] = 0.002484342677554566
$arr[19,0] = 0.002478831348389102
$arr[20,0] = 0.002475366889393245
$arr[21,0] = 0.002477203395437
$arr[22,0] = 0.002484435088424433
$arr[23,0] = 0.002492830271350055
$ws.Range("G2:G25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 1.746282104336743
$arr[1,0] = 1.719032899959458
$arr[2,0] = 1.70324719800395
$arr[3,0] = 1.697051819144448
$arr[4,0] = 1.696037412083456
$arr[5,0] = 1.703162684060757
$arr[6,0] = 1.736690154748274
$arr[7,0] = 1.809958868860932
$arr[8,0] = 1.868410920623177
$arr[9,0] = 1.896014281248483
$arr[10,0] = 1.906613189990537
$arr[11,0] = 1.904324020137452
$arr[12,0] = 1.89688332907707
$arr[13,0] = 1.892344731923515
$arr[14,0] = 1.866627400008156
$arr[15,0] = 1.851110453494584
$arr[16,0] = 1.842280860830556
$arr[17,0] = 1.839307675528005
$arr[18,0] = 1.852752388202575
$arr[19,0] = 1.89906487153209
$arr[20,0] = 1.930184627031252
$arr[21,0] = 1.91349734214549
$arr[22,0] = 1.852009785137966
$arr[23,0] = 1.789328805967443
$ws.Range("I2:I25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.2074473435009168
$arr[1,0] = 0.2067835050896321
$arr[2,0] = 0.2064653220325283
$arr[3,0] = 0.2063581490954505
$arr[4,0] = 0.2063417115617412
$arr[5,0] = 0.2064637855957656
$arr[6,0] = 0.2071998956376149
$arr[7,0] = 0.2093530118713858
$arr[8,0] = 0.2113681746737512
$arr[9,0] = 0.2123791936030059
$arr[10,0] = 0.2127756097624669
$arr[11,0] = 0.2126896310699422
$arr[12,0] = 0.2124115351094105
$arr[13,0] = 0.212242960010343
$arr[14,0] = 0.2113040002660753
$arr[15,0] = 0.2107521364256755
$arr[16,0] = 0.2104435957162707
$arr[17,0] = 0.2103406536367132
$arr[18,0] = 0.2108099645390737
$arr[19,0] = 0.2124928504470347
$arr[20,0] = 0.2136717843445766
$arr[21,0] = 0.2130353288470346
$arr[22,0] = 0.2107837932473444
$arr[23,0] = 0.208694501666173
$ws.Range("J2:J25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.2117037805257098
$arr[1,0] = 0.2115744371489257
$arr[2,0] = 0.2115765548513195
$arr[3,0] = 0.2115979369358953
$arr[4,0] = 0.2116027272131902
$arr[5,0] = 0.2115767601141911
$arr[6,0] = 0.2116422695766076
$arr[7,0] = 0.212417194507708
$arr[8,0] = 0.2133803393411284
$arr[9,0] = 0.2139039825388735
$arr[10,0] = 0.214114560421244
$arr[11,0] = 0.2140686624518224
$arr[12,0] = 0.2139210607405175
$arr[13,0] = 0.2138322500946614
$arr[14,0] = 0.2133478377320159
$arr[15,0] = 0.2130725588641837
$arr[16,0] = 0.2129222743098111
$arr[17,0] = 0.2128727731990097
$arr[18,0] = 0.2131010298467189
$arr[19,0] = 0.213964081588955
$arr[20,0] = 0.2145997413840419
$arr[21,0] = 0.2142539289425258
$arr[22,0] = 0.2130881332683714
$arr[23,0] = 0.2121383810424433
$ws.Range("L2:L25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 2.873581687696401
$arr[1,0] = 2.562605684679681
$arr[2,0] = 2.371325805375818
$arr[3,0] = 2.293303068607429
$arr[4,0] = 2.280343261403573
$arr[5,0] = 2.370273851395496
$arr[6,0] = 2.766433886209825
$arr[7,0] = 3.540180268007646
$arr[8,0] = 4.10623028343673
$arr[9,0] = 4.363110593465422
$arr[10,0] = 4.460285735714251
$arr[11,0] = 4.439361943450422
$arr[12,0] = 4.371107314139522
$arr[13,0] = 4.329286057409945
$arr[14,0] = 4.089429168003562
$arr[15,0] = 3.94211849063862
$arr[16,0] = 3.857331695637754
$arr[17,0] = 3.828614786364199
$arr[18,0] = 3.957806003280837
$arr[19,0] = 4.391158149571254
$arr[20,0] = 4.673791817957863
$arr[21,0] = 4.523002190001307
$arr[22,0] = 3.950713976768498
$arr[23,0] = 3.331249627311138
$ws.Range("N2:N25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 6.773999207379518
$arr[1,0] = 6.656197380233266
$arr[2,0] = 6.587612858722366
$arr[3,0] = 6.560604663501977
$arr[4,0] = 6.556176738047782
$arr[5,0] = 6.58724481047858
$arr[6,0] = 6.732602645548809
$arr[7,0] = 7.047465294632445
$arr[8,0] = 7.297133943381823
$arr[9,0] = 7.414734569018833
$arr[10,0] = 7.459848082973167
$arr[11,0] = 7.450106222564614
$arr[12,0] = 7.41843443562874
$arr[13,0] = 7.399110234349791
$arr[14,0] = 7.289529623125077
$arr[15,0] = 7.223337734753272
$arr[16,0] = 7.185644656866771
$arr[17,0] = 7.172947420900641
$arr[18,0] = 7.230344762001494
$arr[19,0] = 7.427721431513589
$arr[20,0] = 7.560104227502791
$arr[21,0] = 7.489138602731657
$arr[22,0] = 7.227175759008787
$arr[23,0] = 6.959077938266944
$ws.Range("O2:O25").Value = $arr
